$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: quantity (C2) goes from 1 to 2, and the "ligado" boolean (E2) turns on
$ws.Range("C2").Value = 2
$ws.Range("E2").Value = $true

# New row 3: a light bulb ("Lâmpada") object
$ws.Range("A3").Value = "l"
$ws.Range("B3").Value = "Lâmpada"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = $false

# Match default (unstyled) formatting of the new row instead of inheriting column styles
$ws.Range("A3:D3").Style = "Normal"
